$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Status" column (E) for the rows whose task is now finished,
# recording who performed the work (EB or AK). Row 5 ("leancanvas blueprint
# -AK") is marked Complete-AK, all the other newly-completed rows are
# marked Complete-EB. Rows 21-25 previously had no status set at all.
$ws.Range("E2").Value = "Complete-EB"
$ws.Range("E3").Value = "Complete-EB"
$ws.Range("E4").Value = "Complete-EB"
$ws.Range("E5").Value = "Complete-AK"
$ws.Range("E6").Value = "Complete-EB"
$ws.Range("E7").Value = "Complete-EB"
$ws.Range("E8").Value = "Complete-EB"
$ws.Range("E9").Value = "Complete-EB"
$ws.Range("E15").Value = "Complete-EB"
$ws.Range("E16").Value = "Complete-EB"
$ws.Range("E20").Value = "Complete-EB"
$ws.Range("E21").Value = "Complete-EB"
$ws.Range("E22").Value = "Complete-EB"
$ws.Range("E23").Value = "Complete-EB"
$ws.Range("E24").Value = "Complete-EB"
$ws.Range("E25").Value = "Complete-EB"

# Column E now holds longer text ("Complete-EB" / "Complete-AK"), so widen
# it to fit the new content, matching the workbook author's column resize.
$ws.Columns.Item(5).ColumnWidth = 10.8

# The author's cursor ended up on E26 after entering the last status value.
[void]$ws.Range("E26").Select()
